$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.428.95"
$ws.Range("E2").Value = "  -1.21%  "
$ws.Range("D3").Value = "2.511.91"
$ws.Range("E3").Value = "  -0.01%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'572.51"
$ws.Range("E5").Value = "  -0.26%  "
$ws.Range("D6").Value = "'165.92"
$ws.Range("E6").Value = "  -1.25%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").Value = "'0.513"
$ws.Range("E8").Value = "  +0.64%  "
$ws.Range("D9").Value = "2.509.04"
$ws.Range("E9").Value = "  -0.19%  "
$ws.Range("D10").Value = "'0.160"
$ws.Range("E10").Value = "  -1.16%  "
$ws.Range("D11").Value = "'0.168"
$ws.Range("E11").Value = "  -0.30%  "
$ws.Range("D12").Value = "'0.357"
$ws.Range("E12").Value = "  +4.67%  "
$ws.Range("D13").Value = "'4.90"
$ws.Range("E13").Value = "  +1.62%  "
$ws.Range("D14").Value = "2.974.11"
$ws.Range("E14").Value = "  -0.07%  "
$ws.Range("D15").Value = "69.232.13"
$ws.Range("E15").Value = "  -1.35%  "
$ws.Range("D16").Value = "'0.0000175"
$ws.Range("E16").Value = "  -1.72%  "
$ws.Range("D17").Value = "'24.82"
$ws.Range("E17").Value = "  -0.43%  "
$ws.Range("D18").Value = "2.519.92"
$ws.Range("E18").Value = "  -0.14%  "
$ws.Range("D19").Value = "'11.32"
$ws.Range("E19").Value = "  -0.40%  "
$ws.Range("D20").Value = "'7.58"
$ws.Range("E20").Value = "  -0.55%  "
$ws.Range("D21").Value = "'348.41"
$ws.Range("E21").Value = "  -1.77%  "
$ws.Range("D22").Value = "'3.91"
$ws.Range("E22").Value = "  -0.73%  "
$ws.Range("D23").Value = "'2.00"
$ws.Range("E23").Value = "  +2.07%  "
$ws.Range("E24").Value = "  +0.12%  "
$ws.Range("D25").Value = "'70.18"
$ws.Range("E25").Value = "  +1.52%  "
$ws.Range("D26").Value = "'3.93"
$ws.Range("E26").Value = "  -3.20%  "
$ws.Range("D27").Value = "'8.89"
$ws.Range("E27").Value = "  -2.38%  "
$ws.Range("D28").Value = "2.656.61"
$ws.Range("E28").Value = "  +0.19%  "
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("D30").Value = "0.0₃0888"
$ws.Range("E30").Value = "  -2.11%  "
$ws.Range("D31").Value = "'7.83"
$ws.Range("E31").Value = "  +0.07%  "
$ws.Range("D32").Value = "'461.84"
$ws.Range("E32").Value = "  -4.00%  "
$ws.Range("D33").Value = "'1.24"
$ws.Range("E33").Value = "  -4.97%  "
$ws.Range("E34").Value = "  -1.96%  "
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "  -0.11%  "
$ws.Range("B36").Value = "Monero"
$ws.Range("C36").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D36").Value = "'157.28"
$ws.Range("E36").Value = "  +0.58%  "
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").Value = "'0.116"
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("D38").Value = "'19.05"
$ws.Range("E38").Value = "  +0.81%  "
$ws.Range("E39").Value = "  -0.66%  "
$ws.Range("E40").Value = "  -0.05%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D41").Value = "'4.70"
$ws.Range("E41").Value = "  -0.35%  "
$ws.Range("B42").Value = "PolygonEcosystemToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D42").Value = "'0.317"
$ws.Range("E42").Value = "  -0.75%  "
$ws.Range("E43").Value = "  -2.08%  "
$ws.Range("D44").Value = "'38.13"
$ws.Range("E44").Value = "  -0.51%  "
$ws.Range("D45").Value = "'1.12"
$ws.Range("E45").Value = "  -5.97%  "
$ws.Range("D46").Value = "'2.22"
$ws.Range("E46").Value = "  -6.08%  "
$ws.Range("D47").Value = "'142.18"
$ws.Range("E47").Value = "  -0.38%  "
$ws.Range("D48").Value = "'3.47"
$ws.Range("E48").Value = "  -1.38%  "
$ws.Range("D49").Value = "'0.518"
$ws.Range("E49").Value = "  -1.33%  "
$ws.Range("D50").Value = "'0.0730"
$ws.Range("E50").Value = "  +0.14%  "
$ws.Range("D51").Value = "'0.578"
$ws.Range("E51").Value = "  -3.51%  "
